$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 128, shifting the old blank row (128) and summary rows down by one.
$ws.Rows.Item(128).Insert()

# Fill the newly inserted row 128 with the new working-hours entry.
$ws.Cells.Item(128, 1).Value = 2014
$ws.Cells.Item(128, 2).Value = 4
$ws.Cells.Item(128, 3).Value = 15
$ws.Cells.Item(128, 4).Value = 0.88541666666666663
$ws.Cells.Item(128, 5).Value = 0.95833333333333337
$ws.Cells.Item(128, 6).Formula = "=(E128-D128)*24*60"
$ws.Cells.Item(128, 7).Formula = "=F128/60"

# Update selection to reflect where the user ended up editing.
$ws.Range("H128").Select() | Out-Null
